$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 23) mirroring existing row pattern
$ws.Range("A23").Value = Get-Date -Year 2021 -Month 12 -Day 12
$ws.Range("B23").Value = 65

# Update selection to match the diff (C25)
$ws.Range("C25").Select()
